$wb = $excel.ActiveWorkbook

# --- Sheet: Detalle_Pasos (H, I, K columns for rows 2-25) ---
$ws1 = $wb.Worksheets.Item("Detalle_Pasos")

$ws1.Range("H2").Value = 2.539706860597317
$ws1.Range("I2").Value = 1.153522334547124
$ws1.Range("K2").Value = 2.032902704351438
$ws1.Range("H3").Value = 1.723936055715267
$ws1.Range("I3").Value = 2.870361813818601
$ws1.Range("K3").Value = 2.230793648659107
$ws1.Range("H4").Value = 0.9813737097427514
$ws1.Range("I4").Value = 2.365721022510591
$ws1.Range("K4").Value = 1.93969075190823
$ws1.Range("H5").Value = 0.3962622329858431
$ws1.Range("I5").Value = 1.573728926239791
$ws1.Range("K5").Value = 1.646847391502495
$ws1.Range("H6").Value = 0.2704546474462021
$ws1.Range("I6").Value = 0.4218539662615082
$ws1.Range("K6").Value = 1.668739007363945
$ws1.Range("H7").Value = 0.7118746374206125
$ws1.Range("I7").Value = 0.4875529574219293
$ws1.Range("K7").Value = 1.788909442203023
$ws1.Range("H8").Value = 0.5725211483001713
$ws1.Range("I8").Value = 0.5427010461423523
$ws1.Range("K8").Value = 1.700551401900429
$ws1.Range("H9").Value = 0.5872463367462162
$ws1.Range("I9").Value = 0.5562344909419679
$ws1.Range("K9").Value = 1.691681735588474
$ws1.Range("H10").Value = 0.251948138159816
$ws1.Range("I10").Value = 0.7103763938655518
$ws1.Range("K10").Value = 1.6870531869567
$ws1.Range("H11").Value = 0.7250768047332765
$ws1.Range("I11").Value = 0.4224237154221395
$ws1.Range("K11").Value = 1.7749306112228
$ws1.Range("H12").Value = 1.310795264053345
$ws1.Range("I12").Value = 2.00894014431761
$ws1.Range("K12").Value = 1.950960876780078
$ws1.Range("H13").Value = 1.846214690017701
$ws1.Range("I13").Value = 2.430352086146202
$ws1.Range("K13").Value = 1.819000407277076
$ws1.Range("H14").Value = 7.451492629989364
$ws1.Range("I14").Value = 8.660045220928707
$ws1.Range("K14").Value = 5.734071476362839
$ws1.Range("H15").Value = 14.80205645998554
$ws1.Range("I15").Value = 10.02674879476327
$ws1.Range("K15").Value = 18.00500237850628
$ws1.Range("H16").Value = 4.417871956641857
$ws1.Range("I16").Value = 15.47137474156478
$ws1.Range("K16").Value = 8.375835186681101
$ws1.Range("H17").Value = 1.704711103051372
$ws1.Range("I17").Value = 1.218110129946754
$ws1.Range("K17").Value = 3.263949372709502
$ws1.Range("H18").Value = 1.72273913025856
$ws1.Range("I18").Value = 1.311463721485095
$ws1.Range("K18").Value = 2.375325059849129
$ws1.Range("H19").Value = 3.011498615762891
$ws1.Range("I19").Value = 2.146897093238744
$ws1.Range("K19").Value = 2.675805090820443
$ws1.Range("H20").Value = 2.755394513261389
$ws1.Range("I20").Value = 2.086743036905924
$ws1.Range("K20").Value = 2.936485144773834
$ws1.Range("H21").Value = 0.5665686197576254
$ws1.Range("I21").Value = 2.449740384995365
$ws1.Range("K21").Value = 1.869076177355903
$ws1.Range("H22").Value = 2.901349834633621
$ws1.Range("I22").Value = 2.303740687381113
$ws1.Range("K22").Value = 1.395298584892969
$ws1.Range("H23").Value = 1.226069298545284
$ws1.Range("I23").Value = 0.9386866476530393
$ws1.Range("K23").Value = 1.425312474453537
$ws1.Range("H24").Value = 3.261386633004188
$ws1.Range("I24").Value = 2.346841564524471
$ws1.Range("K24").Value = 1.988700028216392
$ws1.Range("H25").Value = 0.9141754155623406
$ws1.Range("I25").Value = 1.30747945547431
$ws1.Range("K25").Value = 1.508925258103195

# --- Sheet: Reliability_Data (Empirical column C) ---
$ws2 = $wb.Worksheets.Item("Reliability_Data")

$ws2.Range("C201").Value = 0.04166666666666666
$ws2.Range("C205").Value = 0.08333333333333333
$ws2.Range("C206").Value = 0.08333333333333333
$ws2.Range("C207").Value = 0.08333333333333333
$ws2.Range("C208").Value = 0.1666666666666667
$ws2.Range("C212").Value = 0.25
$ws2.Range("C222").Value = 0.2916666666666667
$ws2.Range("C223").Value = 0.2916666666666667
$ws2.Range("C231").Value = 0.3333333333333333
$ws2.Range("C232").Value = 0.375
$ws2.Range("C244").Value = 0.4166666666666667
$ws2.Range("C245").Value = 0.4166666666666667
$ws2.Range("C250").Value = 0.5
$ws2.Range("C257").Value = 0.5
$ws2.Range("C271").Value = 0.5833333333333334
$ws2.Range("C272").Value = 0.625
$ws2.Range("C273").Value = 0.625
$ws2.Range("C275").Value = 0.6666666666666666
$ws2.Range("C278").Value = 0.75
$ws2.Range("C279").Value = 0.75
$ws2.Range("C280").Value = 0.7916666666666666
$ws2.Range("C281").Value = 0.8333333333333334
$ws2.Range("C282").Value = 0.8333333333333334
$ws2.Range("C286").Value = 0.9166666666666666
$ws2.Range("C398").Value = 0.125
$ws2.Range("C399").Value = 0.25
$ws2.Range("C419").Value = 0.4583333333333333
$ws2.Range("C420").Value = 0.4583333333333333
$ws2.Range("C421").Value = 0.4583333333333333
$ws2.Range("C433").Value = 0.7083333333333334
$ws2.Range("C448").Value = 0.8333333333333334
$ws2.Range("C449").Value = 0.875
$ws2.Range("C480").Value = 0.9166666666666666
$ws2.Range("C506").Value = 0.0
$ws2.Range("C507").Value = 0.0
$ws2.Range("C509").Value = 0.08333333333333333
$ws2.Range("C511").Value = 0.125
$ws2.Range("C513").Value = 0.25
$ws2.Range("C514").Value = 0.25
$ws2.Range("C517").Value = 0.25
$ws2.Range("C518").Value = 0.25
$ws2.Range("C519").Value = 0.2916666666666667
$ws2.Range("C520").Value = 0.2916666666666667
$ws2.Range("C522").Value = 0.4166666666666667
$ws2.Range("C523").Value = 0.4166666666666667
$ws2.Range("C524").Value = 0.4166666666666667
$ws2.Range("C525").Value = 0.4166666666666667
$ws2.Range("C526").Value = 0.4166666666666667
$ws2.Range("C535").Value = 0.4583333333333333
$ws2.Range("C536").Value = 0.4583333333333333
$ws2.Range("C539").Value = 0.5
$ws2.Range("C542").Value = 0.5833333333333334
$ws2.Range("C543").Value = 0.5833333333333334
$ws2.Range("C544").Value = 0.5833333333333334
$ws2.Range("C545").Value = 0.5833333333333334
$ws2.Range("C546").Value = 0.5833333333333334
$ws2.Range("C557").Value = 0.625
$ws2.Range("C558").Value = 0.625
$ws2.Range("C559").Value = 0.625
$ws2.Range("C560").Value = 0.625
$ws2.Range("C562").Value = 0.6666666666666666
$ws2.Range("C572").Value = 0.6666666666666666
$ws2.Range("C573").Value = 0.6666666666666666
$ws2.Range("C574").Value = 0.6666666666666666
$ws2.Range("C575").Value = 0.6666666666666666
$ws2.Range("C577").Value = 0.7083333333333334
$ws2.Range("C578").Value = 0.7083333333333334
$ws2.Range("C579").Value = 0.7083333333333334
$ws2.Range("C580").Value = 0.7916666666666666
$ws2.Range("C581").Value = 0.7916666666666666
$ws2.Range("C582").Value = 0.7916666666666666
$ws2.Range("C583").Value = 0.7916666666666666
$ws2.Range("C584").Value = 0.7916666666666666
$ws2.Range("C585").Value = 0.8333333333333334
$ws2.Range("C589").Value = 0.9166666666666666
$ws2.Range("C590").Value = 0.9166666666666666
$ws2.Range("C602").Value = 0.0
$ws2.Range("C610").Value = 0.04166666666666666
$ws2.Range("C611").Value = 0.04166666666666666
$ws2.Range("C612").Value = 0.04166666666666666
$ws2.Range("C613").Value = 0.04166666666666666
$ws2.Range("C616").Value = 0.2083333333333333
$ws2.Range("C617").Value = 0.2083333333333333
$ws2.Range("C618").Value = 0.2083333333333333
$ws2.Range("C619").Value = 0.3333333333333333
$ws2.Range("C620").Value = 0.3333333333333333
$ws2.Range("C622").Value = 0.375
$ws2.Range("C623").Value = 0.375
$ws2.Range("C624").Value = 0.375
$ws2.Range("C626").Value = 0.4166666666666667
$ws2.Range("C629").Value = 0.5
$ws2.Range("C631").Value = 0.5416666666666666
$ws2.Range("C632").Value = 0.5416666666666666
$ws2.Range("C633").Value = 0.5416666666666666
$ws2.Range("C634").Value = 0.5416666666666666
$ws2.Range("C635").Value = 0.5416666666666666
$ws2.Range("C636").Value = 0.5416666666666666
$ws2.Range("C637").Value = 0.5416666666666666
$ws2.Range("C638").Value = 0.5416666666666666
$ws2.Range("C639").Value = 0.5833333333333334
$ws2.Range("C640").Value = 0.5833333333333334
$ws2.Range("C641").Value = 0.5833333333333334
$ws2.Range("C642").Value = 0.5833333333333334
$ws2.Range("C643").Value = 0.5833333333333334
$ws2.Range("C655").Value = 0.5833333333333334
$ws2.Range("C656").Value = 0.5833333333333334
$ws2.Range("C657").Value = 0.5833333333333334
$ws2.Range("C658").Value = 0.5833333333333334
$ws2.Range("C659").Value = 0.625
$ws2.Range("C662").Value = 0.75
$ws2.Range("C663").Value = 0.75
$ws2.Range("C681").Value = 0.7916666666666666
$ws2.Range("C682").Value = 0.7916666666666666
$ws2.Range("C683").Value = 0.7916666666666666
$ws2.Range("C684").Value = 0.7916666666666666
$ws2.Range("C685").Value = 0.7916666666666666
$ws2.Range("C686").Value = 0.875
$ws2.Range("C687").Value = 0.875
$ws2.Range("C688").Value = 0.875
$ws2.Range("C689").Value = 0.875
$ws2.Range("C816").Value = 0.08333333333333333
$ws2.Range("C817").Value = 0.125
$ws2.Range("C818").Value = 0.125
$ws2.Range("C819").Value = 0.1666666666666667
$ws2.Range("C820").Value = 0.1666666666666667
$ws2.Range("C821").Value = 0.1666666666666667
$ws2.Range("C822").Value = 0.1666666666666667
$ws2.Range("C823").Value = 0.1666666666666667
$ws2.Range("C824").Value = 0.1666666666666667
$ws2.Range("C825").Value = 0.1666666666666667
$ws2.Range("C826").Value = 0.2916666666666667
$ws2.Range("C827").Value = 0.2916666666666667
$ws2.Range("C828").Value = 0.2916666666666667
$ws2.Range("C829").Value = 0.2916666666666667
$ws2.Range("C830").Value = 0.2916666666666667
$ws2.Range("C831").Value = 0.2916666666666667
$ws2.Range("C832").Value = 0.3333333333333333
$ws2.Range("C833").Value = 0.375
$ws2.Range("C834").Value = 0.375
$ws2.Range("C835").Value = 0.375
$ws2.Range("C836").Value = 0.375
$ws2.Range("C837").Value = 0.375
$ws2.Range("C838").Value = 0.375
$ws2.Range("C839").Value = 0.375
$ws2.Range("C840").Value = 0.375
$ws2.Range("C841").Value = 0.4583333333333333
$ws2.Range("C842").Value = 0.4583333333333333
$ws2.Range("C843").Value = 0.5
$ws2.Range("C844").Value = 0.5
$ws2.Range("C845").Value = 0.5833333333333334
$ws2.Range("C846").Value = 0.5833333333333334
$ws2.Range("C847").Value = 0.6666666666666666
$ws2.Range("C848").Value = 0.6666666666666666
$ws2.Range("C849").Value = 0.7083333333333334
$ws2.Range("C850").Value = 0.7083333333333334
$ws2.Range("C851").Value = 0.75
$ws2.Range("C852").Value = 0.7916666666666666
$ws2.Range("C853").Value = 0.7916666666666666
$ws2.Range("C854").Value = 0.7916666666666666
$ws2.Range("C855").Value = 0.7916666666666666
$ws2.Range("C856").Value = 0.7916666666666666
$ws2.Range("C857").Value = 0.7916666666666666
$ws2.Range("C858").Value = 0.7916666666666666
$ws2.Range("C859").Value = 0.8333333333333334
$ws2.Range("C860").Value = 0.8333333333333334
$ws2.Range("C861").Value = 0.8333333333333334
$ws2.Range("C862").Value = 0.8333333333333334
$ws2.Range("C863").Value = 0.875
$ws2.Range("C864").Value = 0.875
$ws2.Range("C865").Value = 0.9166666666666666
$ws2.Range("C866").Value = 0.9166666666666666
$ws2.Range("C867").Value = 0.9166666666666666
$ws2.Range("C868").Value = 0.9166666666666666
$ws2.Range("C869").Value = 0.9166666666666666
$ws2.Range("C870").Value = 0.9166666666666666
$ws2.Range("C871").Value = 0.9166666666666666
$ws2.Range("C872").Value = 0.9166666666666666
$ws2.Range("C873").Value = 0.9166666666666666
$ws2.Range("C874").Value = 0.9166666666666666
$ws2.Range("C875").Value = 0.9166666666666666
$ws2.Range("C876").Value = 0.9166666666666666
$ws2.Range("C877").Value = 0.9166666666666666
$ws2.Range("C878").Value = 0.9166666666666666
$ws2.Range("C879").Value = 0.9166666666666666
$ws2.Range("C880").Value = 0.9166666666666666
$ws2.Range("C881").Value = 0.9166666666666666
$ws2.Range("C882").Value = 0.9166666666666666
$ws2.Range("C889").Value = 0.9583333333333334
$ws2.Range("C890").Value = 0.9583333333333334
$ws2.Range("C891").Value = 0.9583333333333334
$ws2.Range("C892").Value = 0.9583333333333334
